$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 10: optocoupler part
$ws.Range("B10").Value = "Optokoppler "
$ws.Range("C10").Value = "ACPL-064L-000E"

$ws.Hyperlinks.Add($ws.Range("C10"), "https://www.mouser.de/ProductDetail/Broadcom/ACPL-064L-000E") | Out-Null

$ws.Range("C10").Style = "Link"

$ws.Range("C10").Select()
